$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.528.91'
$ws.Range('E2').Value = '  +0.16%  '
$ws.Range('D3').Value = '1.951.73'
$ws.Range('E3').Value = '  +0.57%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.45'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.622'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.77'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.13%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.377'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0789'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -7.57%  '
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.18'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +6.09%  '
$ws.Range('D13').Value = '2.235.86'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.824'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.44'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.24'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = '1.956.00'
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('D18').Value = '36.430.12'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.12'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.14%  '
$ws.Range('D20').Value = '0.0₃0850'
$ws.Range('E20').Value = '  -1.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '228.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.05'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.41%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.70%  '
$ws.Range('E25').Value = '  +2.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.143'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.69%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.57%  '
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.31%  '
$ws.Range('E30').Value = '  +19.36%  '
$ws.Range('E31').Value = '  +1.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.76'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.54%  '
$ws.Range('E33').Value = '  -0.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.43'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.19%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.27'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +4.20%  '
$ws.Range('E37').Value = '  +9.16%  '
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -12.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0967'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.78%  '
$ws.Range('E41').Value = '  -0.86%  '
$ws.Range('E42').Value = '  +1.15%  '
$ws.Range('E43').Value = '  +0.11%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.71'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '1.359.82'
$ws.Range('E45').Value = '  +1.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.45'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.98%  '
$ws.Range('E47').Value = '  -0.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.11'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.83'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.44%  '
$ws.Range('D51').Value = '2.131.82'
$ws.Range('E51').Value = '  +0.72%  '
